$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# birth_date and hire_date columns change type from VARCHAR(30) to DATE
$ws.Range("B10").Value = "DATE"
$ws.Range("B14").Value = "DATE"

# Select the D column's data range (selection moved from G14 to D1:D34)
$ws.Range("D1:D34").Select()
